$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, old C->D, etc.)
$ws.Columns("B:B").Insert()

# The newly inserted column should take on the same width as the old column B (40 chars).
# ColumnWidth round-trips through a pixel conversion in this engine, so compensate
# by the same constant offset (5/7) that a width of 40 comes back as 40 + 5/7.
$ws.Columns("B:B").ColumnWidth = 39.285714285714285

# Fill in the new "button_text" column with header rows + button labels
$ws.Range("B1").Value = "button_text"
$ws.Range("B2").Value = "버튼 텍스트"
$ws.Range("B3").Value = "string"
$ws.Range("B4").Value = "공격하기"
$ws.Range("B5").Value = "빈틈 노리기"
$ws.Range("B6").Value = "공격 회피하기"
$ws.Range("B7").Value = "저지하기"

# Move selection to where the user was last working
[void]$ws.Range("B8").Select()
